# "stop the train correctly red line"
# Move each station's Beacon (column E = flag, F/G = B0/B1 payload) off the
# "Switch #" row and onto the station-name row directly above it, widen the
# B0/B1 columns (F:G) to fit the new content, and scroll the frozen pane
# down to where the edits were made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 7 / SHADYSIDE: beacon info moves from row 9 (E/F) to row 8 (E/F) ---
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "SHADYSIDE,TRUE,BOTH,TRUE"
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()

# --- Block 21 / SWISSVILLE: beacon info moves from row 23 (E/G) to row 22 (E/G) ---
$ws.Range("E22").Value = 1
$ws.Range("G22").Value = "SWISSVILLE,TRUE,BOTH,FALSE"
$ws.Range("E23").ClearContents()
$ws.Range("G23").ClearContents()

# --- Block 45 / FIRST AVE: beacon info moves from row 47 (E/G) to row 46 (E/G) ---
$ws.Range("E46").Value = 1
$ws.Range("G46").Value = "FIRST AVE,TRUE,BOTH,FALSE"
$ws.Range("E47").ClearContents()
$ws.Range("G47").ClearContents()

# --- Block 47 / STATION SQUARE: beacon info moves from row 48 (E/F) to row 49 (E/F) ---
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = "STATION SQUARE,TRUE,BOTH,FALSE"
$ws.Range("E48").ClearContents()
$ws.Range("F48").ClearContents()

# --- Widen the B0/B1 beacon columns now that they hold longer strings ---
$ws.Range("F1:G1").ColumnWidth = 39.36328125

# --- Scroll the frozen view down to where the edits were made ---
$win = $excel.ActiveWindow
$win.ScrollRow = 12
$win.ScrollColumn = 1
$ws.Range("E17").Select()
